$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1514.6666
$ws.Range("I33").Value = 938.8889
$ws.Range("K33").Value = 938.8889
$ws.Range("M33").Value = -709.8889

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3040
$ws.Range("I62").Value = 2200
$ws.Range("K62").Value = 2200
$ws.Range("M62").Value = -1576

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3040
$ws.Range("I65").Value = 2200
$ws.Range("K65").Value = 11000
$ws.Range("M65").Value = -7880

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1591.4286
$ws.Range("I98").Value = 1108.8889
$ws.Range("K98").Value = 1108.8889
$ws.Range("M98").Value = 389.1111000000001

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2499.0908
$ws.Range("I100").Value = 2436.25
$ws.Range("K100").Value = 2436.25
$ws.Range("M100").Value = -1895.25

# ALC row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1497.5
$ws.Range("I121").Value = 530
$ws.Range("J121").Value = 1650.2632
$ws.Range("K121").Value = 1590
$ws.Range("L121").Value = 4950.7896
$ws.Range("M121").Value = 157
$ws.Range("N121").Value = -8444.7896

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1591.4286
$ws.Range("I122").Value = 1108.8889
$ws.Range("K122").Value = 3326.6667
$ws.Range("M122").Value = -876.6666999999998

# ALC row 123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""

# ALC row 126
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H126").Value = 55486.668
$ws.Range("J126").Value = 55486.668
$ws.Range("L126").Value = 55486.668
$ws.Range("N126").Value = -65366.668

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1696.2307
$ws.Range("I135").Value = 1637.7273
$ws.Range("J135").Value = 2018
$ws.Range("K135").Value = 14739.5457
$ws.Range("L135").Value = 18162
$ws.Range("M135").Value = -12204.5457
$ws.Range("N135").Value = -23232

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1321.2222
$ws.Range("I137").Value = 790.2857
$ws.Range("J137").Value = 1507.05
$ws.Range("K137").Value = 2370.8571
$ws.Range("L137").Value = 4521.15
$ws.Range("M137").Value = 179.1428999999998
$ws.Range("N137").Value = -9621.15

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4797.5884
$ws.Range("I138").Value = 4486.091
$ws.Range("K138").Value = 13458.273
$ws.Range("M138").Value = -8318.273000000001

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 452959.44
$ws.Range("I32").Value = 593336.9
$ws.Range("J32").Value = 13110.2
$ws.Range("K32").Value = 593336.9
$ws.Range("L32").Value = 13110.2
$ws.Range("M32").Value = -593049.9
$ws.Range("N32").Value = -13684.2

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2192
$ws.Range("I45").Value = 1600
$ws.Range("J45").Value = 3714.2856
$ws.Range("K45").Value = 1600
$ws.Range("L45").Value = 3714.2856
$ws.Range("M45").Value = -1223
$ws.Range("N45").Value = -4468.2856

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4139.091
$ws.Range("I61").Value = 4500
$ws.Range("K61").Value = 4500
$ws.Range("M61").Value = -4288

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1605.0526
$ws.Range("J74").Value = 1460.3334
$ws.Range("L74").Value = 1460.3334
$ws.Range("N74").Value = -3208.3334

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1605.0526
$ws.Range("J77").Value = 1460.3334
$ws.Range("L77").Value = 7301.666999999999
$ws.Range("N77").Value = -16037.667

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1176.9445
$ws.Range("I122").Value = 793.7778
$ws.Range("J122").Value = 1560.1111
$ws.Range("K122").Value = 2381.3334
$ws.Range("L122").Value = 4680.3333
$ws.Range("M122").Value = 68.66660000000002
$ws.Range("N122").Value = -9580.3333

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3056.7415
$ws.Range("I132").Value = 2356.587
$ws.Range("J132").Value = 5740.6665
$ws.Range("K132").Value = 7069.761
$ws.Range("L132").Value = 17221.9995
$ws.Range("M132").Value = -4539.761
$ws.Range("N132").Value = -22281.9995

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4139.091
$ws.Range("I136").Value = 4500
$ws.Range("K136").Value = 13500
$ws.Range("M136").Value = -10950

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 648.7037
$ws.Range("I94").Value = 639.5833
$ws.Range("J94").Value = 656
$ws.Range("K94").Value = 639.5833
$ws.Range("L94").Value = 656
$ws.Range("M94").Value = -188.5833
$ws.Range("N94").Value = -1558

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3498.5938
$ws.Range("I134").Value = 3963.6667
$ws.Range("J134").Value = 2900.6428
$ws.Range("K134").Value = 11891.0001
$ws.Range("L134").Value = 8701.928400000001
$ws.Range("M134").Value = -9356.000100000001
$ws.Range("N134").Value = -13771.9284

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1454.2972
$ws.Range("I31").Value = 1509.05
$ws.Range("J31").Value = 1389.8823
$ws.Range("K31").Value = 1509.05
$ws.Range("L31").Value = 1389.8823
$ws.Range("M31").Value = -1214.05
$ws.Range("N31").Value = -1979.8823

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1454.2972
$ws.Range("I34").Value = 1509.05
$ws.Range("J34").Value = 1389.8823
$ws.Range("K34").Value = 1509.05
$ws.Range("L34").Value = 1389.8823
$ws.Range("M34").Value = -1307.05
$ws.Range("N34").Value = -1793.8823

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 13891499
$ws.Range("I132").Value = 2574.5715
$ws.Range("J132").Value = 33335994
$ws.Range("K132").Value = 7723.7145
$ws.Range("L132").Value = 100007982
$ws.Range("M132").Value = -5193.7145
$ws.Range("N132").Value = -100013042

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1406.3572
$ws.Range("I134").Value = 997.5
$ws.Range("J134").Value = 2428.5
$ws.Range("K134").Value = 2992.5
$ws.Range("L134").Value = 7285.5
$ws.Range("M134").Value = -457.5
$ws.Range("N134").Value = -12355.5

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 851.0175
$ws.Range("I68").Value = 621.55554
$ws.Range("J68").Value = 956.9231
$ws.Range("K68").Value = 1864.66662
$ws.Range("L68").Value = 2870.7693
$ws.Range("M68").Value = -1053.66662
$ws.Range("N68").Value = -4492.7693

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 851.0175
$ws.Range("I71").Value = 621.55554
$ws.Range("J71").Value = 956.9231
$ws.Range("K71").Value = 5593.99986
$ws.Range("L71").Value = 8612.3079
$ws.Range("M71").Value = -1537.99986
$ws.Range("N71").Value = -16724.3079

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1181.8379
$ws.Range("I107").Value = 213.4
$ws.Range("J107").Value = 1540.5186
$ws.Range("K107").Value = 640.2
$ws.Range("L107").Value = 4621.5558
$ws.Range("M107").Value = 1279.8
$ws.Range("N107").Value = -8461.5558

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 17130
$ws.Range("I133").Value = 780
$ws.Range("J133").Value = 20400
$ws.Range("K133").Value = 2340
$ws.Range("L133").Value = 61200
$ws.Range("M133").Value = 2720
$ws.Range("N133").Value = -71320

# CUL row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2383.3572
$ws.Range("I136").Value = 2033.3636
$ws.Range("J136").Value = 3666.6667
$ws.Range("K136").Value = 6100.0908
$ws.Range("L136").Value = 11000.0001
$ws.Range("M136").Value = -1000.0908
$ws.Range("N136").Value = -21200.0001

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3779.8667
$ws.Range("I132").Value = 3712.5
$ws.Range("J132").Value = 3856.8572
$ws.Range("K132").Value = 11137.5
$ws.Range("L132").Value = 11570.5716
$ws.Range("M132").Value = -8607.5
$ws.Range("N132").Value = -16630.5716

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2000.5714
$ws.Range("I7").Value = 1900.8
$ws.Range("J7").Value = 2250
$ws.Range("K7").Value = 1900.8
$ws.Range("L7").Value = 2250
$ws.Range("M7").Value = -1788.8
$ws.Range("N7").Value = -2474

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 8929.5
$ws.Range("I93").Value = 11844.556
$ws.Range("J93").Value = 3682.4
$ws.Range("K93").Value = 11844.556
$ws.Range("L93").Value = 3682.4
$ws.Range("M93").Value = -10596.556
$ws.Range("N93").Value = -6178.4

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3059.875
$ws.Range("I100").Value = 2447
$ws.Range("J100").Value = 3672.75
$ws.Range("K100").Value = 2447
$ws.Range("L100").Value = 3672.75
$ws.Range("M100").Value = -1906
$ws.Range("N100").Value = -4754.75

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3421.6
$ws.Range("I122").Value = 3702.6667
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 11108.0001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -8658.000100000001
$ws.Range("N122").Value = -13900

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2000.5714
$ws.Range("I126").Value = 1900.8
$ws.Range("J126").Value = 2250
$ws.Range("K126").Value = 5702.4
$ws.Range("L126").Value = 6750
$ws.Range("M126").Value = -3232.4
$ws.Range("N126").Value = -11690

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3843.6453
$ws.Range("I132").Value = 3334.4736
$ws.Range("J132").Value = 4649.8335
$ws.Range("K132").Value = 10003.4208
$ws.Range("L132").Value = 13949.5005
$ws.Range("M132").Value = -7473.4208
$ws.Range("N132").Value = -19009.5005

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2316.5
$ws.Range("I136").Value = 2600
$ws.Range("J136").Value = 899
$ws.Range("K136").Value = 7800
$ws.Range("L136").Value = 2697
$ws.Range("M136").Value = -5250
$ws.Range("N136").Value = -7797

# WVR row 9
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 70007
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = ""

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4121.0293
$ws.Range("I96").Value = 2442.625
$ws.Range("J96").Value = 4637.4614
$ws.Range("K96").Value = 2442.625
$ws.Range("L96").Value = 4637.4614
$ws.Range("M96").Value = -1069.625
$ws.Range("N96").Value = -7383.4614

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 883.5
$ws.Range("I107").Value = 594
$ws.Range("J107").Value = 1007.5714
$ws.Range("K107").Value = 1782
$ws.Range("L107").Value = 3022.7142
$ws.Range("M107").Value = 138
$ws.Range("N107").Value = -6862.7142

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10421326
$ws.Range("I132").Value = 5295.5835
$ws.Range("J132").Value = 41669416
$ws.Range("K132").Value = 15886.7505
$ws.Range("L132").Value = 125008248
$ws.Range("M132").Value = -13356.7505
$ws.Range("N132").Value = -125013308
